$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.806.48"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.968.91"
$ws.Range("E3").Value = "  +4.22%  "
$ws.Range("E4").Value = "  -1.22%  "
$ws.Range("D5").Value = "'251.58"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").Value = "'0.7135"
$ws.Range("E6").Value = "  +51.60%  "
$ws.Range("D7").Value = "'0.9889"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").Value = "'0.3237"
$ws.Range("E8").Value = "  +11.75%  "
$ws.Range("D9").Value = "'25.81"
$ws.Range("E9").Value = "  +16.31%  "
$ws.Range("D10").Value = "'0.06843"
$ws.Range("E10").Value = "  +5.39%  "
$ws.Range("D11").Value = "'0.8318"
$ws.Range("E11").Value = "  +14.67%  "
$ws.Range("D12").Value = "'101.45"
$ws.Range("E12").Value = "  +6.05%  "
$ws.Range("D13").Value = "'0.07929"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "1.947.18"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").Value = "'5.376"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").Value = "'276.36"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "30.793.44"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'13.84"
$ws.Range("E18").Value = "  +6.01%  "
$ws.Range("D19").Value = "'0.000007711"
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").Value = "'5.609"
$ws.Range("E20").Value = "  +6.18%  "
$ws.Range("D21").Value = "2.201.77"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "'0.9916"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "'0.9885"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "'6.642"
$ws.Range("E24").Value = "  +6.08%  "
$ws.Range("D25").Value = "'9.512"
$ws.Range("E25").Value = "  +4.88%  "
$ws.Range("D26").Value = "'164.42"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'19.48"
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("E28").Value = "  +31.87%  "
$ws.Range("D29").Value = "'2.161"
$ws.Range("E29").Value = "  +14.02%  "
$ws.Range("D30").Value = "'1.349"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "'1.546"
$ws.Range("E31").Value = "  +5.16%  "
$ws.Range("D32").Value = "'4.489"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").Value = "'4.334"
$ws.Range("E33").Value = "  +4.61%  "
$ws.Range("D34").Value = "'0.05034"
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("D35").Value = "'1.202"
$ws.Range("E35").Value = "  +6.77%  "
$ws.Range("D36").Value = "'0.7301"
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'0.01975"
$ws.Range("D39").Value = "'2.915"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").Value = "'6.548"
$ws.Range("E40").Value = "  +5.42%  "
$ws.Range("D41").Value = "'77.62"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "'0.4612"
$ws.Range("E42").Value = "  +8.70%  "
$ws.Range("D43").Value = "'2.052"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").Value = "'0.8391"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").Value = "'0.9894"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").Value = "'10.01"
$ws.Range("E46").Value = "  +4.58%  "
$ws.Range("D47").Value = "'102.12"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").Value = "'7.356"
$ws.Range("E48").Value = "  +5.74%  "
$ws.Range("D49").Value = "'36.01"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "'0.4202"
$ws.Range("E50").Value = "  +6.84%  "
$ws.Range("D51").Value = "'936.69"
$ws.Range("E51").Value = "  +2.75%  "
